$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$new.Name = "Demand Charge"

$new.Range("B1").Value = "direct PV use"
$new.Range("A2").Value = "LA"
$new.Range("B2").Value = 6211.2916123429659
$new.Range("A3").Value = "Boulder"
$new.Range("B3").Value = 7007.1402678343975
$new.Range("B2:B3").NumberFormat = "0"
$new.Columns.Item(2).ColumnWidth = 11.83203125

foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
